$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8321464493434974
$ws.Range("C2").Value = 0.2183703710160785
$ws.Range("E2").Value = 0.4255945681936879
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.002368596919982519
$ws.Range("N2").Value = 0.809833198744407
$ws.Range("O2").Value = 1.385518127233553

$ws.Range("B3").Value = 0.729054412652431
$ws.Range("C3").Value = 0.1940060492390501
$ws.Range("E3").Value = 0.3711622844637077
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.002371697075888117
$ws.Range("N3").Value = 0.8139744630372263
$ws.Range("O3").Value = 1.369872005639763

$ws.Range("B4").Value = 0.6656786958320424
$ws.Range("C4").Value = 0.1789658641554013
$ws.Range("E4").Value = 0.3378381819432263
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.002373699025721346
$ws.Range("N4").Value = 0.8169254432623134
$ws.Range("O4").Value = 1.361729970833522

$ws.Range("B5").Value = 0.6398339292775574
$ws.Range("C5").Value = 0.1728168069009826
$ws.Range("E5").Value = 0.3242808676558582
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.002374539672497145
$ws.Range("N5").Value = 0.8182307626882164
$ws.Range("O5").Value = 1.358777953036594

$ws.Range("B6").Value = 0.6355413190343313
$ws.Range("C6").Value = 0.1717945537360777
$ws.Range("E6").Value = 0.3220309884541592
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.002374680763735815
$ws.Range("N6").Value = 0.8184537207387024
$ws.Range("O6").Value = 1.358309803947321

$ws.Range("B7").Value = 0.6653302188868793
$ws.Range("C7").Value = 0.1788830168290474
$ws.Range("E7").Value = 0.3376552546630052
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.002373710262232629
$ws.Range("N7").Value = 0.8169426309832488
$ws.Range("O7").Value = 1.361688680442683

$ws.Range("B8").Value = 0.7966166152938854
$ws.Range("C8").Value = 0.2099863239570539
$ws.Range("E8").Value = 0.4068048501122945
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.002369645471884339
$ws.Range("N8").Value = 0.8111764472638185
$ws.Range("O8").Value = 1.379817883164463

$ws.Range("B9").Value = 1.053443613538832
$ws.Range("C9").Value = 0.2703380557669561
$ws.Range("E9").Value = 0.5432776199618417
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.002362451756905457
$ws.Range("N9").Value = 0.8031037082158008
$ws.Range("O9").Value = 1.427103674656166

$ws.Range("B10").Value = 1.241751292722427
$ws.Range("C10").Value = 0.3142864214171937
$ws.Range("E10").Value = 0.6442182070340152
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.002357635125946625
$ws.Range("N10").Value = 0.7991400336670296
$ws.Range("O10").Value = 1.469160004652849

$ws.Range("B11").Value = 1.327335126915386
$ws.Range("C11").Value = 0.3341946887317135
$ws.Range("E11").Value = 0.6903167798395629
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.002355544536823321
$ws.Range("N11").Value = 0.7977632915569899
$ws.Range("O11").Value = 1.489916409744865

$ws.Range("B12").Value = 1.359731904967589
$ws.Range("C12").Value = 0.3417212611233538
$ws.Range("E12").Value = 0.7078014320347421
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.00235476725399721
$ws.Range("N12").Value = 0.7973032027205278
$ws.Range("O12").Value = 1.498012735531404

$ws.Range("B13").Value = 1.352755222309895
$ws.Range("C13").Value = 0.3401008268330656
$ws.Range("E13").Value = 0.7040345130603356
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.002354934017460243
$ws.Range("N13").Value = 0.7973995676180721
$ws.Range("O13").Value = 1.496258493319175

$ws.Range("B14").Value = 1.330000673091661
$ws.Range("C14").Value = 0.3348141509577545
$ws.Range("E14").Value = 0.6917546763031481
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.002355480301618804
$ws.Range("N14").Value = 0.7977242124510013
$ws.Range("O14").Value = 1.490577747580545

$ws.Range("B15").Value = 1.316061276044934
$ws.Range("C15").Value = 0.3315743082675056
$ws.Range("E15").Value = 0.6842366502823864
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.002355816786700345
$ws.Range("N15").Value = 0.7979310422830679
$ws.Range("O15").Value = 1.48712898335026

$ws.Range("B16").Value = 1.236156552084822
$ws.Range("C16").Value = 0.3129836632921297
$ws.Range("E16").Value = 0.6412093712781797
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.002357773767044234
$ws.Range("N16").Value = 0.799238581394377
$ws.Range("O16").Value = 1.46783644850035

$ws.Range("B17").Value = 1.187117037111364
$ws.Range("C17").Value = 0.3015572531278679
$ws.Range("E17").Value = 0.6148612462464627
$ws.Range("F17").Value = 0.6400460337215605
$ws.Range("G17").Value = 0.002359000002124388
$ws.Range("N17").Value = 0.8001498716084825
$ws.Range("O17").Value = 1.456419151574124

$ws.Range("B18").Value = 1.158903479073388
$ws.Range("C18").Value = 0.2949771814681412
$ws.Range("E18").Value = 0.5997233284111161
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.002359714766091054
$ws.Range("N18").Value = 0.800714154873404
$ws.Range("O18").Value = 1.450004993858613

$ws.Range("B19").Value = 1.149349616601512
$ws.Range("C19").Value = 0.2927479301395408
$ws.Range("E19").Value = 0.5946007126186572
$ws.Range("F19").Value = 0.6191636801734006
$ws.Range("G19").Value = 0.002359958401085531
$ws.Range("N19").Value = 0.8009121062266615
$ws.Range("O19").Value = 1.447859426986867

$ws.Range("B20").Value = 1.192338140572474
$ws.Range("C20").Value = 0.3027744331224937
$ws.Range("E20").Value = 0.6176642916540942
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.002358868488255307
$ws.Range("N20").Value = 0.8000487100129874
$ws.Range("O20").Value = 1.457618711896032

$ws.Range("B21").Value = 1.336684566002532
$ws.Range("C21").Value = 0.3363673097689457
$ws.Range("E21").Value = 0.6953607810138749
$ws.Range("F21").Value = 0.7228739723492197
$ws.Range("G21").Value = 0.002355319455010072
$ws.Range("N21").Value = 0.7976271943955737
$ws.Range("O21").Value = 1.492239885010491

$ws.Range("B22").Value = 1.430953254092969
$ws.Range("C22").Value = 0.3582506816890998
$ws.Range("E22").Value = 0.7463048973279029
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.002353083726647291
$ws.Range("N22").Value = 0.796401606001055
$ws.Range("O22").Value = 1.516245567910232

$ws.Range("B23").Value = 1.380646920336062
$ws.Range("C23").Value = 0.3465777124260399
$ws.Range("E23").Value = 0.7190992311472257
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.002354269337325344
$ws.Range("N23").Value = 0.7970230753472123
$ws.Range("O23").Value = 1.503306236256748

$ws.Range("B24").Value = 1.189977742215547
$ws.Range("C24").Value = 0.3022241799144467
$ws.Range("E24").Value = 0.616397003956962
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.002358927915041994
$ws.Range("N24").Value = 0.8000943194142565
$ws.Range("O24").Value = 1.457075924269873

$ws.Range("B25").Value = 0.9840322917638105
$ws.Range("C25").Value = 0.2540801152023562
$ws.Range("E25").Value = 0.5062501904638026
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.002364315177292259
$ws.Range("N25").Value = 0.8049418045789238
$ws.Range("O25").Value = 1.413038689045976

